$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B1 (styled numeric 0) becomes A1 - use Copy so the cell style (border/bold/
# centered) travels with it, same as Excel's own "Copy" command would do.
$ws.Range("B1").Copy($ws.Range("A1"))

# B2:B7 (inline-string item names) become A2:A7, replacing the old numeric
# index values that used to live there. Strip the old numeric-column style
# first so the moved-in text ends up with the default (no explicit) format,
# then write over the values.
$ws.Range("A2:A7").ClearFormats()
$ws.Range("A2").Value = $ws.Range("B2").Value2
$ws.Range("A3").Value = $ws.Range("B3").Value2
$ws.Range("A4").Value = $ws.Range("B4").Value2
$ws.Range("A5").Value = $ws.Range("B5").Value2
$ws.Range("A6").Value = $ws.Range("B6").Value2
$ws.Range("A7").Value = $ws.Range("B7").Value2

# Column B is no longer used - clear it out entirely.
$ws.Range("B1:B7").Clear()
